$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reference rows to append: CaB6 (rows 35-37) and Fe4C (rows 38-40)
$data = @(
    @("CaB6", 0.535, 74, 5, 80),
    @("CaB6", 0.486, 17, -3, 15),
    @("CaB6", 0.43, 8, -21, -12),
    @("Fe4C", 0.535, 112, -1, 111),
    @("Fe4C", 0.486, 26, 58, 84),
    @("Fe4C", 0.43, -32, 108, 77)
)

$startRow = 35
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

$excel.ActiveWindow.ScrollRow = 5
$ws.Range("A40").Select()
